$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows that keep the same region (rows 2-6)
$ws.Range("B2").Value = 83
$ws.Range("C2").Value = 78

$ws.Range("B3").Value = 42
$ws.Range("C3").Value = 38

$ws.Range("B4").Value = 19
$ws.Range("C4").Value = 18

$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 15

$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 12

# Rows 7-12 get re-ordered regions (Marlborough moves up, others shift) plus
# new values for some rows.
$ws.Range("A7").Value = "Marlborough"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 1

$ws.Range("A8").Value = "Nelson"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 0

$ws.Range("A9").Value = "Taranaki"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = "Bay of Plenty"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 0

$ws.Range("A11").Value = "Hawke$([char]0x2019)s Bay"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Manawatu-Whanganui"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0

# Row 14: South Canterbury -> TBC, with new counts
$ws.Range("A14").Value = "TBC"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 2
